$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 0.757176135198727
$ws.Range("P2").Value = 0.757176135198727
$ws.Range("S2").Value = 0.757176135198727
$ws.Range("T2").Value = 0.757176135198727

# Row 3 updates
$ws.Range("M3").Value = 0.6217929999999999
$ws.Range("N3").Value = 1.865379
$ws.Range("O3").Value = 0.242823864801273
$ws.Range("P3").Value = 0.2428238648012731
$ws.Range("Q3").Value = 1.099928188866
$ws.Range("R3").Value = 9.899353699794
$ws.Range("S3").Value = 0.242823864801273
$ws.Range("T3").Value = 0.2428238648012731
